# Update column F values on sheet "展览" and sheet "全部类型"
$wb = $excel.ActiveWorkbook

# Mapping of row -> new value for column F
$updates = @{
    2  = 628
    3  = 2193
    4  = 83
    5  = 13017
    6  = 71
    7  = 115
    8  = 514
    10 = 1173
    11 = 975
    12 = 13741
    13 = 14286
    22 = 1086
    25 = 5355
    26 = 933
    27 = 16
    28 = 296
    29 = 8
    30 = 11
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
